$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2457912457912458
$ws.Range("C2").Value = 0.4646464646464646
$ws.Range("J2").Value = 0.0101010101010101
$ws.Range("P2").Value = 0.1683501683501684
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("B3").Value = 0.03973509933774835
$ws.Range("C3").Value = 0.05298013245033113
$ws.Range("J3").Value = 0.02649006622516556
$ws.Range("P3").Value = 0.6887417218543046
$ws.Range("S3").Value = 0.1920529801324503
$ws.Range("P4").Value = 0.8148148148148148
$ws.Range("S4").Value = 0.1851851851851852
$ws.Range("B6").Value = 0.06796116504854369
$ws.Range("D6").Value = 0.01456310679611651
$ws.Range("F6").Value = 0.08737864077669903
$ws.Range("J6").Value = 0.2524271844660194
$ws.Range("O6").Value = 0.009708737864077669
$ws.Range("Q6").Value = 0.0970873786407767
$ws.Range("R6").Value = 0.06796116504854369
$ws.Range("S6").Value = 0.4029126213592233
$ws.Range("B7").Value = 0.1229050279329609
$ws.Range("D7").Value = 0.0223463687150838
$ws.Range("F7").Value = 0.0782122905027933
$ws.Range("J7").Value = 0.1284916201117318
$ws.Range("O7").Value = 0.0111731843575419
$ws.Range("Q7").Value = 0.1564245810055866
$ws.Range("R7").Value = 0.0782122905027933
$ws.Range("S7").Value = 0.4022346368715084
$ws.Range("B8").Value = 0.09115281501340483
$ws.Range("D8").Value = 0.008042895442359249
$ws.Range("E8").Value = 0.005361930294906166
$ws.Range("F8").Value = 0.05361930294906166
$ws.Range("J8").Value = 0.1394101876675603
$ws.Range("O8").Value = 0.01876675603217158
$ws.Range("Q8").Value = 0.1876675603217158
$ws.Range("R8").Value = 0.0938337801608579
$ws.Range("S8").Value = 0.4021447721179625
$ws.Range("B9").Value = 0.1005917159763314
$ws.Range("D9").Value = 0.01183431952662722
$ws.Range("F9").Value = 0.05325443786982249
$ws.Range("J9").Value = 0.1301775147928994
$ws.Range("O9").Value = 0.03550295857988166
$ws.Range("Q9").Value = 0.2011834319526627
$ws.Range("R9").Value = 0.05917159763313609
$ws.Range("S9").Value = 0.408284023668639
$ws.Range("B10").Value = 0.1109185441941074
$ws.Range("D10").Value = 0.0147313691507799
$ws.Range("E10").Value = 0.0008665511265164644
$ws.Range("F10").Value = 0.07192374350086655
$ws.Range("J10").Value = 0.1360485268630849
$ws.Range("O10").Value = 0.01039861351819757
$ws.Range("Q10").Value = 0.1949740034662045
$ws.Range("R10").Value = 0.0684575389948007
$ws.Range("S10").Value = 0.391681109185442
$ws.Range("G11").Value = 0.1366906474820144
$ws.Range("J11").Value = 0.07194244604316546
$ws.Range("K11").Value = 0.2014388489208633
$ws.Range("L11").Value = 0.5611510791366906
$ws.Range("S11").Value = 0.02877697841726619
$ws.Range("G12").Value = 0.7034883720930233
$ws.Range("J12").Value = 0.1395348837209302
$ws.Range("K12").Value = 0.01744186046511628
$ws.Range("L12").Value = 0.06976744186046512
$ws.Range("S12").Value = 0.06976744186046512
$ws.Range("G13").Value = 0.7647058823529411
$ws.Range("J13").Value = 0.2058823529411765
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("F15").Value = 0.03883495145631068
$ws.Range("H15").Value = 0.145631067961165
$ws.Range("I15").Value = 0.07766990291262135
$ws.Range("J15").Value = 0.3980582524271845
$ws.Range("K15").Value = 0.07766990291262135
$ws.Range("M15").Value = 0.01456310679611651
$ws.Range("O15").Value = 0.06310679611650485
$ws.Range("S15").Value = 0.1844660194174757
$ws.Range("F16").Value = 0.01734104046242774
$ws.Range("H16").Value = 0.1734104046242775
$ws.Range("I16").Value = 0.04046242774566474
$ws.Range("J16").Value = 0.4566473988439306
$ws.Range("K16").Value = 0.08670520231213873
$ws.Range("M16").Value = 0.02312138728323699
$ws.Range("O16").Value = 0.05202312138728324
$ws.Range("S16").Value = 0.1502890173410405
$ws.Range("F17").Value = 0.02688172043010753
$ws.Range("H17").Value = 0.1720430107526882
$ws.Range("I17").Value = 0.06989247311827956
$ws.Range("J17").Value = 0.3924731182795699
$ws.Range("K17").Value = 0.1129032258064516
$ws.Range("M17").Value = 0.01881720430107527
$ws.Range("N17").Value = 0.002688172043010753
$ws.Range("O17").Value = 0.08602150537634409
$ws.Range("S17").Value = 0.1182795698924731
$ws.Range("F18").Value = 0.01324503311258278
$ws.Range("H18").Value = 0.1258278145695364
$ws.Range("I18").Value = 0.1258278145695364
$ws.Range("J18").Value = 0.4172185430463576
$ws.Range("K18").Value = 0.09933774834437085
$ws.Range("M18").Value = 0.01324503311258278
$ws.Range("O18").Value = 0.1324503311258278
$ws.Range("S18").Value = 0.0728476821192053
$ws.Range("F19").Value = 0.01964133219470538
$ws.Range("H19").Value = 0.197267292912041
$ws.Range("I19").Value = 0.08710503842869342
$ws.Range("J19").Value = 0.372331340734415
$ws.Range("K19").Value = 0.1067463706233988
$ws.Range("M19").Value = 0.01707941929974381
$ws.Range("N19").Value = 0.003415883859948762
$ws.Range("O19").Value = 0.06660973526900085
$ws.Range("S19").Value = 0.1298035866780529
